# Update concise_ms csv pattern: refresh the "Marking" and "Total" rows
# on the quiz marksheet (rows 11-12) to the new scoring values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right marks 4 -> 5, Wrong marks -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 ("Total"): Right total 68 -> 85, Wrong total 0 -> -0, summary label updated
$ws.Range("B12").Value = 85
$ws.Range("C12").Value = -0
$ws.Range("E12").Value = "85.0/140"
